$d = $word.ActiveDocument

# Paragraph 2: "Gijs de Vries, s1854526" -> set language to Dutch (nl-NL)
$p2 = $d.Paragraphs(2)
$p2.Range.LanguageID = 1043

# Paragraph 3: "Revision 0.1" -> set language to Dutch (nl-NL)
$p3 = $d.Paragraphs(3)
$p3.Range.LanguageID = 1043
